# Reverse the order of comma-separated "Recorded By" names (column G) for
# every row where the value currently starts with "System, " (i.e. "System"
# is listed first). This flips entries like:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"   -> "backup@backdoor.com, system, System"
# Rows whose value does not start with "System, " (already reordered, a
# single value, or doesn't contain "System" at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val

    if ($text.StartsWith("System, ")) {
        $parts = $text -split ", "
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}
